$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 'crop athletic leggings'
$ws.Range("A2").Value = 'kid basketball knee pads'
$ws.Range("A3").Value = 'men''s spandex leggings'
$ws.Range("A4").Value = 'hex pants'
$ws.Range("A5").Value = 'medical compression pants'
$ws.Range("A6").Value = 'fitness tights for men'
$ws.Range("A7").Value = 'the rock mens basketball'
$ws.Range("A8").Value = 'black men tights'
$ws.Range("A9").Value = 'd man basketball'
$ws.Range("A10").Value = 'asics compression pants'
$ws.Range("A11").Value = 'knee armor knee pads'
$ws.Range("A12").Value = 'adidas tights for men'
$ws.Range("A13").Value = 'slide on knee pads'
$ws.Range("A14").Value = 'men athletic tights'
$ws.Range("A15").Value = 'blue mens compression pants'
$ws.Range("A16").Value = 'sport tights for men'
$ws.Range("A17").Value = 'compression basketball tights'
$ws.Range("A18").Value = 'kids compression pants'
$ws.Range("A19").Value = 'usa tights men'
$ws.Range("A20").Value = 'razor knee pads'
$ws.Range("A21").Value = 'knee pads addidas'
$ws.Range("A22").Value = 'knee pads hunting'
$ws.Range("A23").Value = 'knee pads leggings'
$ws.Range("A24").Value = 'knee pads elbow pads youth'
$ws.Range("A25").Value = 'knee pads slim'
$ws.Range("A26").Value = 'basketball pants adidas'
$ws.Range("A27").Value = 'mens tights navy'
$ws.Range("A28").Value = 'mens adidas basketball pants'
$ws.Range("A29").Value = 'compression pants men 3xl'
$ws.Range("A30").Value = 'compression pants 2xu'
$ws.Range("A31").Value = 'pro x knee pad'
$ws.Range("A32").Value = 'puma compression pants men'
$ws.Range("A33").Value = 'men''s basketball pants'
$ws.Range("A34").Value = 'gray baseball pants youth girls'
$ws.Range("A35").Value = 'mens basketball jacket'
$ws.Range("A36").Value = 'photography knee pads'
$ws.Range("A37").Value = 'men''s tights leggings'
$ws.Range("A38").Value = 'baseball pants men grey'
$ws.Range("A39").Value = 'venom compression pants'
$ws.Range("A40").Value = 'padded knee tights'
$ws.Range("A41").Value = 'ua compression pants'
$ws.Range("A42").Value = 'men workout tights'
$ws.Range("A43").Value = 'flag compression pants'
$ws.Range("A44").Value = 'cool knee pads'
$ws.Range("A45").Value = 'navy compression leggings'
$ws.Range("A46").Value = 'force knee pads'
$ws.Range("A47").Value = 'mens wrestling pants'
$ws.Range("A48").Value = 'mens capri compression pants'
$ws.Range("A49").Value = 'woman compression pants'
$ws.Range("A50").Value = 'purple knee pads'
$ws.Range("A51").Value = 'reebok knee pads'
$ws.Range("A52").Value = 'venum compression pants men'
$ws.Range("A53").Value = 'purple athletic leggings'
$ws.Range("A54").Value = 'thermal compression pants'
$ws.Range("A55").Value = 'addidas knee pads'
$ws.Range("A56").Value = 'jordan mens tights'
$ws.Range("A57").Value = 'purple compression pants men'
$ws.Range("A58").Value = 'russell compression pants'
$ws.Range("A59").Value = 'blue knee pads for basketball'
$ws.Range("A60").Value = 'elbow knee pad'
$ws.Range("A61").Value = 'tights mens'
$ws.Range("A62").Value = 'super compression leggings'
$ws.Range("A63").Value = 'mens leggins'
$ws.Range("A64").Value = 'knee pads for teens'
$ws.Range("A65").Value = 'green mens compression pants'
$ws.Range("A66").Value = 'pants with padded knees'
$ws.Range("A67").Value = 'compression with pads'
$ws.Range("A68").Value = 'knee pads for sleeping'
$ws.Range("A69").Value = 'mens winter compression pants'
$ws.Range("A70").Value = 'knee pads for hvac'
$ws.Range("A71").Value = 'yoga pants with knee pads'
$ws.Range("A72").Value = 'black pants with knee pads'
$ws.Range("A73").Value = 'kids compression knee pads'
$ws.Range("A74").Value = 'elite basketball pants'
$ws.Range("A75").Value = 'nike leggings mens'
$ws.Range("A76").Value = 'compression pants men 3 pack'
$ws.Range("A77").Value = 'compression pants baseball'
$ws.Range("A78").Value = 'colored compression leggings'
$ws.Range("A79").Value = 'mens compression pants with pockets'
$ws.Range("A80").Value = 'mens compression pants xxl'
$ws.Range("A81").Value = 'knee pads wheels'
$ws.Range("A82").Value = 'ua basketball knee pads'
$ws.Range("A83").Value = 'protective knee pad'
$ws.Range("A84").Value = 'grey knee pads basketball'
$ws.Range("A85").Value = 'orange knee pads for basketball'
$ws.Range("A86").Value = 'adidas youth compression pants'
$ws.Range("A87").Value = 'copper compression tights for men'
$ws.Range("A88").Value = 'basketball knee pads youth boys mcdavid'
$ws.Range("A89").Value = 'nike youth basketball knee pads'
$ws.Range("A90").Value = 'nike kneepads'
$ws.Range("A91").Value = 'knee pad strap'
$ws.Range("A92").Value = 'internal knee pad'
$ws.Range("A93").Value = 'basketball legings'
$ws.Range("A94").Value = 'basketball knee sleves'
$ws.Range("A95").Value = 'under armour baseball pants men'
$ws.Range("A96").Value = 'cold gear compression leggings men'
$ws.Range("A97").Value = 'mcgregor knee pads'
$ws.Range("A98").Value = 'knee compression pants men'
$ws.Range("A99").Value = 'compression knee pads for basketball for kids'
$ws.Range("A100").Value = 'tesla compression pant'
